$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.542.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.617.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.52"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.613.55"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.632"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.08"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.96%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.72"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.197.96"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.96"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.615.89"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.500.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -7.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.34"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.10"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.46"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.59"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.24"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.12"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0817"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.25"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +20.19%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.245.44"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.07"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0447"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.75%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.138"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.13%  "
